$d = $word.ActiveDocument

# The "Recorded Classes" table is the first (only) table in the document.
$t = $d.Tables.Item(1)

# Append a new row at the bottom of the table (after the "22nd June" row).
$newRow = $t.Rows.Add()

# --- Column 1: Date -> "24th June" with "th" superscripted ---
$cell1 = $newRow.Cells.Item(1)
$cell1.Range.Text = "24th June"
$dateRange = $cell1.Range
$supRange = $d.Range($dateRange.Start + 2, $dateRange.Start + 4)
$supRange.Font.Superscript = $true

# --- Column 2: Topics -> "Prototype" ---
$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "Prototype"

# --- Column 3: YouTube Link -> hyperlink to the new recording ---
$cell3 = $newRow.Cells.Item(3)
$url = "https://www.youtube.com/watch?v=9MjMVv2GnBY"
$cell3.Range.Text = $url
$linkTextRange = $d.Range($cell3.Range.Start, $cell3.Range.Start + $url.Length)
$d.Hyperlinks.Add($linkTextRange, $url) | Out-Null
